$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 updates (California - Los Angeles) per July 06, 2020 05:39:58 AM CT run
$ws.Range("B16").Value = 44016
$ws.Range("C16").Value = 114993
$ws.Range("D16").Value = 3487
$ws.Range("E16").Value = 3112
$ws.Range("F16").Value = 366
$ws.Range("G16").Value = 4.67
$ws.Range("H16").Value = 11.3
$ws.Range("K16").Value = 66661
$ws.Range("L16").Value = 3240
